$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Data updates - Training Dashboard, row 3 (LOTO (SOPs) training line)
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training Dashboard")

# PERIOD TO EXPIRE: 92 -> 84
$ws.Range("H3").Value = 84

# LAST UPDATE: 08-Sep-2025 -> 16-Sep-2025 (stored as plain text, like the
# rest of that column - the leading apostrophe stops it being read back as
# a real date serial)
$ws.Range("I3").Value = "'16-Sep-2025"

# ------------------------------------------------------------------
# Header formatting - make the dark-blue title/column-header bars use
# white bold text (both sheets) instead of plain black bold text.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

$ws.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
